# Scen_NoBndBiomass.xlsx edit
# - rename Sheet1 -> NoBiomassBound
# - add a new "readme" sheet after it (becomes the active tab)
# - populate readme with a purpose / data-source / user-modifiable summary

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "NoBiomassBound"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "readme"

# Base font for the whole readme sheet
$ws2.Range("A1:A13").Font.Name = "Calibri"
$ws2.Range("A1:A13").Font.Size = 11

# Content - headers first, then their values (keeps shared-string order tidy)
$ws2.Range("A1").Value = "Workbook purpose:"
$ws2.Range("A5").Value = "Data source:"
$ws2.Range("A12").Value = "User-modifiable?"

$ws2.Range("A2").Value = 'Sets an "infinite" biomass supply'
$ws2.Range("A6").Value = "NA"
$ws2.Range("A13").Value = "Yes"

# Bold section headers
$ws2.Range("A1").Font.Bold = $true
$ws2.Range("A5").Font.Bold = $true
$ws2.Range("A12").Font.Bold = $true

$ws2.Activate()
$ws2.Range("A14").Select() | Out-Null
